$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.338.41'
$ws.Range('E2').Value = '  -4.36%  '

$ws.Range('D3').Value = '2.943.07'
$ws.Range('E3').Value = '  -1.06%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '555.54'
$ws.Range('E5').Value = '  -2.85%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.71'
$ws.Range('E6').Value = '  +4.31%  '

$ws.Range('E7').Value = '  +0.12%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.511'
$ws.Range('E8').Value = '  +2.18%  '

$ws.Range('D9').Value = '2.936.77'
$ws.Range('E9').Value = '  -1.19%  '

$ws.Range('E10').Value = '  -3.75%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '4.77'
$ws.Range('E11').Value = '  -6.06%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.443'
$ws.Range('E12').Value = '  +1.30%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000220'
$ws.Range('E13').Value = '  -0.98%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.48'
$ws.Range('E14').Value = '  -0.52%  '

$ws.Range('E15').Value = '  +1.31%  '

$ws.Range('D16').Value = '3.420.12'
$ws.Range('E16').Value = '  -0.73%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.75'
$ws.Range('E17').Value = '  +9.21%  '

$ws.Range('D18').Value = '2.934.19'
$ws.Range('E18').Value = '  -0.50%  '

$ws.Range('D19').Value = '57.332.96'
$ws.Range('E19').Value = '  -4.33%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '415.37'
$ws.Range('E20').Value = '  -4.24%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.07'
$ws.Range('E21').Value = '  +0.01%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.679'
$ws.Range('E22').Value = '  +2.83%  '

$ws.Range('E23').Value = '  -0.76%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.89'
$ws.Range('E24').Value = '  +0.65%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '78.92'
$ws.Range('E25').Value = '  +0.02%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.997'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.49'
$ws.Range('E28').Value = '  -1.22%  '

$ws.Range('E29').Value = '  +4.00%  '

$ws.Range('E30').Value = '  +5.05%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.11'
$ws.Range('E31').Value = '  -0.97%  '

$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '24.99'
$ws.Range('E32').Value = '  -1.24%  '

$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.102'
$ws.Range('E33').Value = '  +9.67%  '

$ws.Range('E34').Value = '  +0.30%  '

$ws.Range('B35').Value = 'Mantle'
$ws.Range('C35').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.931'
$ws.Range('E35').Value = '  -1.59%  '

$ws.Range('B36').Value = 'Stacks'
$ws.Range('C36').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.09'
$ws.Range('E36').Value = '  -3.79%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '48.45'
$ws.Range('E37').Value = '  -2.31%  '

$ws.Range('D38').Value = '0.0₃0675'
$ws.Range('E38').Value = '  +3.00%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.38'
$ws.Range('E39').Value = '  +6.12%  '

$ws.Range('E40').Value = '  +3.53%  '

$ws.Range('E41').Value = '  -2.97%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.108'
$ws.Range('E42').Value = '  -1.14%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '375.79'
$ws.Range('E43').Value = '  -1.20%  '

$ws.Range('D44').Value = '2.630.99'
$ws.Range('E44').Value = '  +0.46%  '

$ws.Range('E45').Value = '  -0.05%  '

$ws.Range('E46').Value = '  +0.86%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.76'
$ws.Range('E47').Value = '  +2.68%  '

$ws.Range('E48').Value = '  +2.23%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.97'
$ws.Range('E49').Value = '  -0.12%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.26'
$ws.Range('E50').Value = '  -0.34%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.99'
$ws.Range('E51').Value = '  +0.30%  '
